$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# Insert two new rows before row 6 (pushes everything from row 6 down by 2)
$ws.Rows.Item(6).Resize(2).Insert()

$ws.Range("A6").Value = "Prior distribution for fluxes (uniform or normal)"
$ws.Range("B6").Value = "normal"
$ws.Range("A7").Value = "Prior distribution for thermodynamic quantities (uniform or normal)"
$ws.Range("B7").Value = "normal"

$ws.Range("A6:A7").Select()
